$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sammer Manuel block: AA30 status changed from 1 to 0.2 ---
$ws.Range("AA30").Value = 0.2

# --- Row 31: new "Backend Security" entry added to Sammer Manuel block ---
$ws.Range("W31").Value = "Backend Security"
$ws.Range("X31").Value = 42816
$ws.Range("Y31").Value = 0.59027777777777779
$ws.Range("Z31").Value = 0.65972222222222221
$ws.Range("AA31").Value = 0.4

# --- Row 36: new entries added to Lamprecht Daniel and Ruhdorfer Alexander blocks ---
$ws.Range("C36").Value = "User Management"
$ws.Range("D36").Value = 42816
$ws.Range("E36").Value = 0.59027777777777779
$ws.Range("F36").Value = 0.65972222222222221
$ws.Range("G36").Value = 0.5

$ws.Range("M36").Value = "Google Maps routing angefangen einzubauen"
$ws.Range("N36").Value = 42816
$ws.Range("O36").Value = 0.59027777777777779
$ws.Range("P36").Value = 0.65972222222222221
$ws.Range("Q36").Value = 0.2

# --- View state: move selection to AA33 (topLeftCell scroll target is not
# persisted by this runtime's writer, but the selection is) ---
$ws.Range("AA33").Select()
